# Add Victor Grazi's email address to the title slide's author textbox,
# and resize/reposition the textbox to fit the extra line.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)   # "TextBox 4"

# Reposition (points; COM uses points, OOXML stores EMU - 1 pt = 12700 EMU).
$shp.Left  = 41.585
$shp.Top   = 440.41512
$shp.Width = 230.3031

# Append a new paragraph with the email address, matching the existing run formatting.
$tr = $shp.TextFrame.TextRange
$tr.InsertAfter([char]13 + "vgrazi@gmail.com") | Out-Null

# Set the final height after the text/autofit change is applied.
$shp.Height = 94.5141
